# Generate Report for Handback
#
# - "Status" text changes from "Ready for handoff" to
#   "Handed back: in sync with en-US" everywhere it appears (Overview!E2,
#   Overview!F2, zh-cn!C2, de-de!C2 all share this value).
# - zh-cn / de-de sheets: fill in "Latest Target File" (I2) and
#   "Latest Handback File" (J2) with the handed-back file names, add a
#   hyperlink on I2 (same target as the existing A2 hyperlink), and stamp
#   "Latest Handback DateTime" (K2) with the handback timestamp.
# - Widen a few columns that now hold longer file-name/date values.

$wb = $excel.ActiveWorkbook

$mdName  = "9db4f506-09c0-41b9-b5f5-b7ddff47b76a.md"
$mdUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f364f6003c6ee15ff7d8e3199207c0783deec7/e2e/9db4f506-09c0-41b9-b5f5-b7ddff47b76a.md"

# ---- Overview sheet --------------------------------------------------
# "Status" text ("Ready for handoff" -> "Handed back: in sync with en-US")
# is shared by Overview!E2/F2 and the per-language Status cells (C2) on the
# zh-cn / de-de sheets, since they all point at the same shared string.
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---- zh-cn sheet -------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("I2").Value = $mdName
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, "", "", $mdName)
# Match the workbook's existing custom "HyperLink" cell style (used by A2)
# instead of the auto-inserted theme-coloured built-in hyperlink style.
$zh.Range("I2").Font.Color = 15570276
$zh.Range("I2").Font.Underline = 2
$zh.Range("J2").Value = "9db4f506-09c0-41b9-b5f5-b7ddff47b76a.f97d646ee18c8449b3c72bfcb12479a71bb09980.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-28 11:07:26"
$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# ---- de-de sheet ---------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("I2").Value = $mdName
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, "", "", $mdName)
$de.Range("I2").Font.Color = 15570276
$de.Range("I2").Font.Underline = 2
$de.Range("J2").Value = "9db4f506-09c0-41b9-b5f5-b7ddff47b76a.f97d646ee18c8449b3c72bfcb12479a71bb09980.de-de.xlf"
$de.Range("K2").Value = "2016-08-28 11:07:33"
$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
